$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.309.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.682.04"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.679.72"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -5.54%  "
$ws.Range("E10").Value = "  -8.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.22"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -8.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.303.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -9.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.681.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.356.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -8.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.650"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -8.22%  "
$ws.Range("E24").Value = "  -4.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.828.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000126"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -11.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.07"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -10.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -9.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.89%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.65%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.652.63"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.161"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -10.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.13"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -6.58%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0903"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.64%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -6.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "164.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.87"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -14.91%  "
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000272"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -10.56%  "
